# Updated cryptos list on Sat Oct 14 19:09:24 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for the crypto
# table, and fixes the ShibaInu / BitcoinCash row ordering (rows 18-19
# were swapped upstream).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, [string]$value)
    # Columns B/C/E never look like bare numbers, but D frequently does
    # ("1.01", "207.53", ...). Excel's COM layer auto-coerces a bare
    # numeric-looking string into a real number (dropping meaningful
    # trailing zeros, e.g. "217.30" -> 217.3), so force text formatting
    # first, write the value, then restore the cell style so no stray
    # formatting is left behind on cells that should stay plain.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "26.986.54"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.560.69"
$ws.Range("E3").Value = "  +0.57%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.20%  "

# Row 5 - BNB
Set-TextCell "D5" "207.53"
$ws.Range("E5").Value = "  +0.41%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.72%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.14%  "

# Row 8 - Solana
$ws.Range("E8").Value = "  +1.85%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.26%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.01%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.14%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextCell "D12" "1.783.25"
$ws.Range("E12").Value = "  +0.56%  "

# Row 13 - WrappedEther
Set-TextCell "D13" "1.546.14"
$ws.Range("E13").Value = "  -0.33%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.73%  "

# Row 15 - Polygon
Set-TextCell "D15" "0.520"
$ws.Range("E15").Value = "  +0.98%  "

# Row 16 - Litecoin
Set-TextCell "D16" "62.04"
$ws.Range("E16").Value = "  +0.59%  "

# Row 17 - WrappedBTC
Set-TextCell "D17" "26.989.51"
$ws.Range("E17").Value = "  +0.35%  "

# Rows 18/19 - ShibaInu and BitcoinCash swapped order upstream
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell "D18" "217.30"
$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D19" "0.0₃0705"
$ws.Range("E19").Value = "  +2.47%  "

# Row 20 - Chainlink
Set-TextCell "D20" "7.38"
$ws.Range("E20").Value = "  +2.30%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.42%  "

# Row 23 - Avalanche
Set-TextCell "D23" "9.22"
$ws.Range("E23").Value = "  +0.05%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -1.42%  "

# Row 25 - Monero
Set-TextCell "D25" "153.61"

# Row 26 - Cosmos
$ws.Range("E26").Value = "  +0.40%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  +1.29%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +1.51%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  -0.09%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +0.76%  "

# Row 31 - PancakeSwap
Set-TextCell "D31" "1.12"
$ws.Range("E31").Value = "  +2.08%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.75%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextCell "D33" "3.12"
$ws.Range("E33").Value = "  +3.90%  "

# Row 34 - Maker
Set-TextCell "D34" "1.423.44"

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +3.04%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  +9.55%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  +1.34%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.82%  "

# Row 39 - ImmutableX
Set-TextCell "D39" "0.532"
$ws.Range("E39").Value = "  +2.08%  "

# Row 40 - ARBITRUM
Set-TextCell "D40" "0.809"
$ws.Range("E40").Value = "  +0.00%  "

# Row 41 - PaxDollar
Set-TextCell "D41" "1.01"
$ws.Range("E41").Value = "  -0.11%  "

# Row 42 - FraxShare
Set-TextCell "D42" "5.71"
$ws.Range("E42").Value = "  -0.12%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  +2.91%  "

# Row 44 - WEMIXToken
$ws.Range("E44").Value = "  +1.77%  "

# Row 45 - Aave
Set-TextCell "D45" "64.91"
$ws.Range("E45").Value = "  +1.95%  "

# Row 46 - RenderToken
$ws.Range("E46").Value = "  -0.05%  "

# Row 47 - RocketPoolETH
Set-TextCell "D47" "1.696.18"
$ws.Range("E47").Value = "  +0.51%  "

# Row 48 - Quant
Set-TextCell "D48" "87.42"
$ws.Range("E48").Value = "  +1.41%  "

# Row 49 - Cronos
Set-TextCell "D49" "0.0523"
$ws.Range("E49").Value = "  +0.19%  "

# Row 50 - BabyDogeCoin
Set-TextCell "D50" "0.0₆0100"
$ws.Range("E50").Value = "  +2.39%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  -0.22%  "
